# Apply the "add eda and change order" edit to the socks product sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price column (D) keeps storing its values as plain text
# (the original sheet stored "49.95 zł" as text; the new values like
# "49.95" look numeric, so force a text number format first to stop
# Excel from auto-converting them to real numbers).
$ws.Range("D2:D13").NumberFormat = "@"

# Header row: rename price -> price_in_PLN
$ws.Range("D1").Value = "price_in_PLN"

# Row 2: HT6538 / black -> price stripped of currency suffix
$ws.Range("D2").Value = "49.95"

# Row 3: IB7814 / navy blue
$ws.Range("D3").Value = "49.95"

# Row 4: IB7817 / red
$ws.Range("D4").Value = "49.95"

# Row 5: IB7818 / royal blue -> blue
$ws.Range("C5").Value = "blue"
$ws.Range("D5").Value = "49.95"

# Row 6: IB7819 / green
$ws.Range("D6").Value = "49.95"

# Row 7: IB7820 / maroon -> burgundy
$ws.Range("C7").Value = "burgundy"
$ws.Range("D7").Value = "49.95"

# Rows 8-10: reorder the IB7821/IB7822/IB7823 entries
$ws.Range("B8").Value = "IB7821"
$ws.Range("C8").Value = "orange"
$ws.Range("D8").Value = "49.95"

$ws.Range("B9").Value = "IB7822"
$ws.Range("C9").Value = "blue"
$ws.Range("D9").Value = "49.95"

$ws.Range("B10").Value = "IB7823"
$ws.Range("C10").Value = "mint"
$ws.Range("D10").Value = "49.95"

# Row 11: IB7815 / yellow
$ws.Range("D11").Value = "49.95"

# Row 12: IB7816 / black -> grey
$ws.Range("C12").Value = "grey"
$ws.Range("D12").Value = "49.95"

# Row 13: new record - socks / IB7813 / white / 49.95
$ws.Range("A13").Value = "socks"
$ws.Range("B13").Value = "IB7813"
$ws.Range("C13").Value = "white"
$ws.Range("D13").Value = "49.95"
